# Rename the worksheet "dati" -> "data" to reflect the synced folder
# naming (matlabfilesuk/ <-> matlabMfilesUK/, capPlots -> ch8Plots).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dati")
$ws.Name = "data"
